$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# The phrase ", tìm kiếm" used to be split across two runs with a leftover
# "_GoBack" bookmark sitting between them (", tìm k" | bookmark | "iếm").
# Re-typing/merging it back into a single run also drops the stale
# bookmark, matching the target XML.
$d.Content.Find.Execute(", tìm kiếm", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", tìm kiếm", 2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# "Mô tả ngắn gọn nội dung nhiệm vụ 4:" appears twice in the document (task
# 4's own description, and task 5's description which was mislabeled as
# "4"). Only the second one - the one describing the renter-management task
# ("người thuê" / "số nhà thuê" ...) - needs to become "... nhiệm vụ 5:".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("nhiệm vụ 4") -and $t.Contains("người thuê")) {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Find.Execute("Mô tả ngắn gọn nội dung nhiệm vụ 4:", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "Mô tả ngắn gọn nội dung nhiệm vụ 5:", 2) | Out-Null
}

# --- Change 3 -------------------------------------------------------------
# Fix the typo "chưa chả" -> "chưa trả" in that same paragraph's sentence.
$d.Content.Find.Execute("số nợ chưa chả", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "số nợ chưa trả", 2) | Out-Null
